# Update column B ("Price") values on Sheet1 to reflect new randomized
# coin prices, now tracked with higher (10 decimal place) accuracy for
# coins that have very low prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$updates = @{
    2  = 512
    3  = 1839
    4  = 1659
    5  = 693
    6  = 870
    7  = 1637
    8  = 495
    9  = 52
    10 = 463
    11 = 1027
    12 = 1274
    13 = 372
    14 = 2544
    15 = 2087
    16 = 2638
    17 = 1586
    18 = 2469
    19 = 2137
    20 = 1521
    21 = 2682
    22 = 2900
    23 = 1230
    24 = 1087
    26 = 3437
    27 = 2405
    28 = 1808
    29 = 131
    30 = 109
    31 = 3873
    32 = 3704
    33 = 2502
    34 = 5268
    35 = 1698
    36 = 1975
    37 = 5994
    38 = 66
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
